# empty_1sheet_template.xlsx edit:
#   - A1 gets the value 0
#   - selection/active cell moves on to B1 (as after typing into A1 and hitting Enter)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0
$ws.Range("B1").Select() | Out-Null
